$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Finalizacao calc e adicao na massa"
# The closing record that used to sit way down at row 10 (REG-401 / 12000000)
# is the real end of the "massa" block, so it moves up to row 4 - right after
# the existing records - and the old, stray underline-formatted row 4
# placeholder cells go away.

# 1) Strip the leftover underline formatting from the destination row before
#    putting any content in it, so the whole A4:H4 row ends up as plain,
#    non-underlined cells (this also clears what used to be on C4/H4).
$ws.Range("A4:H4").Font.Underline = 2
$ws.Range("A4:H4").Font.Underline = -4142

# 2) Bring the row 10 values up into row 4. Quote-prefix them on entry so they
#    stay stored as text (they read as plain numbers otherwise), matching how
#    they were stored before.
$ws.Range("A4").Value = "'REG-401"
$ws.Range("H4").Value = "'12000000"

# 3) The source row (10) is now fully vacated - delete it outright so the
#    sheet's used range/dimension shrinks back down to row 4.
$ws.Range("A10:H10").EntireRow.Delete() | Out-Null

# 4) Leave the selection where the edit landed.
$ws.Range("A4:H4").Select() | Out-Null
